$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row to the table (grows Table2 from A1:D6 to A1:D7)
$tbl = $ws.ListObjects.Item("Table2")
$newRow = $tbl.ListRows.Add()

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Phiri"
$ws.Range("C7").Value = 54

# Carry over the "Hyperlink" cell format used by the other EMAIL cells
# before writing the value, so the new cell matches the existing column style.
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = "phiri@gmail,com"

# Add the mailto hyperlink for the new email cell
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:phiri@gmail,com")

# Re-apply the column's cell format, since adding the hyperlink resets it
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selected cell shown in the sheet view
$ws.Range("G15").Select()
